# "Add a tile and makefile" - update the Sokoban tile-editor grid on the
# "Visual" sheet with a new tile pattern, and draw a border box around the
# editable 8x8 grid (A1:H8). The "Binary" and "Code" sheets are driven by
# formulas referencing "Visual", so they recompute automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Visual")

# New tile bitmap (1 = filled, 0 = empty) for rows 1..8, columns A..H.
$grid = @(
  @(0,0,0,1,1,1,0,0),
  @(0,0,0,1,0,1,0,0),
  @(0,0,0,1,1,0,0,0),
  @(0,1,1,1,1,1,1,0),
  @(1,0,1,1,1,0,0,1),
  @(0,0,1,1,1,1,0,0),
  @(0,1,0,0,0,0,1,0),
  @(0,1,1,0,0,0,1,1)
)

for ($r = 1; $r -le 8; $r++) {
    $row = $grid[$r - 1]
    for ($c = 1; $c -le 8; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        if ($row[$c - 1] -eq 1) {
            $cell.Value = 1
        } else {
            $cell.ClearContents()
        }
    }
}

# Draw a thin box border around the whole tile grid.
$ws.Range("A1:H8").BorderAround(1)

Write-Output "tile grid updated"
